$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.288.11"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "2.554.36"
$ws.Range("E3").Value = "  +4.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.23"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -1.14%  "

$ws.Range("D9").Value = "2.552.51"
$ws.Range("E9").Value = "  +4.64%  "

$ws.Range("E10").Value = "  +0.24%  "

$ws.Range("E11").Value = "  -2.07%  "

$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("E13").Value = "  +0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.30"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.48%  "

$ws.Range("D15").Value = "3.010.37"
$ws.Range("E15").Value = "  +4.79%  "

$ws.Range("D16").Value = "63.192.29"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "2.484.86"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("E19").Value = "  +2.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "334.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.01"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.55%  "

$ws.Range("E26").Value = "  +5.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.48"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +11.66%  "

$ws.Range("E29").Value = "  +3.37%  "

$ws.Range("E30").Value = "  +7.05%  "

$ws.Range("E31").Value = "  +4.08%  "

$ws.Range("E32").Value = "  +1.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "177.67"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.77%  "

$ws.Range("E34").Value = "  +6.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "415.71"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +11.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.398"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.89"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.82%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.18%  "

$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.75"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.23"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.76"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.56%  "

$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("E46").Value = "  +1.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0966"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.56%  "

$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0236"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.43"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.91%  "

$ws.Range("E51").Value = "  +3.94%  "
